$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the list of tasks ("Quoi") in column H (used as the source list for the
# data-validation dropdown in column E). They must be created in this order
# so the shared-strings table ends up with the same ordering as the target.
$ws.Range("H4").Value = "Analyse et état de l'art"
$ws.Range("H5").Value = "Réalisation du modèle"
$ws.Range("H6").Value = "Réalisation de l'application "
$ws.Range("H7").Value = "Tests et validations"
$ws.Range("H8").Value = "Gestion du projet, documentation et présentation"

# The old "Cahier des charges" entries are replaced by the new task
# "Gestion du projet, documentation et présentation" (this removes the last
# reference to "Cahier des charges" so it drops out of the shared strings).
$ws.Range("E2").Value = "Gestion du projet, documentation et présentation"
$ws.Range("E3").Value = "Gestion du projet, documentation et présentation"
$ws.Range("E4").Value = "Gestion du projet, documentation et présentation"
$ws.Range("E5").Value = "Gestion du projet, documentation et présentation"

# New work journal entries.
$ws.Range("C5").Value = 0.41944444444444445

$ws.Range("A6").Value = 45474
$ws.Range("B6").Value = 0.4201388888888889
$ws.Range("C6").Value = 0.70833333333333337

# Widen column E a bit and give column H (the hidden helper list) a width.
$ws.Columns.Item(5).ColumnWidth = 45.17
$ws.Columns.Item(8).ColumnWidth = 44.5

# Turn E1:E1048576 into a dropdown list validated against the H4:H8 list.
$ws.Range("E1:E1048576").Validation.Add(3, 1, 1, "=`$H`$4:`$H`$8") | Out-Null

# Update selection to reflect where the user ended up working.
$ws.Range("C7").Select() | Out-Null
